$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-15 (columns A-F)
$data = @(
    @(45952, 7351.84865254281, 6815.44354432553, 9620, 3189.52974202526, 16.0405535979494),
    @(45953, 7132.87012029834, 6631.57894710238, 3620, 3268.41568801476, 261.666443129881),
    @(45954, 7309.07075906196, 6349.00977144898, 3620, 3407.07375204529, 255.670146812261),
    @(45955, 3049.34035190743, 4847.97569749922, 3620, 3591.87642559605, 200.827171795636),
    @(45956, 2982.57345832622, 4912.73604312149, 3620, 3621.72192057632, 204.769081820742),
    @(45957, 7224.82758554897, 7617.19907167211, 3620, 3510.48975102868, 312.820367612533),
    @(45958, 7224.82758554897, 7825.85177481706, 3620, 3510.48975102868, 321.514230243573),
    @(45959, 7224.82758554897, 7662.27340143543, 3620, 3510.48975102868, 314.698464686005),
    @(45960, 7224.82758554897, 7359.10564219196, 3620, 3510.48975102868, 302.066474717526),
    @(45961, 6107.35191671871, 5897.54092068411, 3620, 3445.87206897926, 238.475541235974),
    @(45962, 2557.16655997904, 4035.97663082571, 9004, 3794.11896080932, -48.912683681874),
    @(45963, 2557.16655997904, 4006.88026636202, 9004, 3794.11896080932, -50.1250322011942),
    @(45964, 8571.71877057272, 7082.94609113525, 9004, 4682.66778538349, 115.067244854947),
    @(45965, 8571.71877057272, 7393.5310507481, 9004, 4682.66778538349, 128.008284838816)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
}

$wb.Save()
